# Append two new daily rows (2026-01-21 / Excel serial 46043) for both
# charging stations at the bottom of the existing data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42: 四方坪站 (Sifangping station)
$ws.Cells.Item(42, 1).Value = 46043
$ws.Cells.Item(42, 2).Value = "四方坪站"
$ws.Cells.Item(42, 3).Value = 16705.51
$ws.Cells.Item(42, 4).Value = 13522.3
$ws.Cells.Item(42, 5).Value = 5438.17
$ws.Cells.Item(42, 6).Value = 670

# Row 43: 高岭站 (Gaoling station)
$ws.Cells.Item(43, 1).Value = 46043
$ws.Cells.Item(43, 2).Value = "高岭站"
$ws.Cells.Item(43, 3).Value = 5705.48
$ws.Cells.Item(43, 4).Value = 4821.66
$ws.Cells.Item(43, 5).Value = 1578.7
$ws.Cells.Item(43, 6).Value = 238

# Move the visible selection down to the newly added rows, as in the source.
$ws.Range("I42").Select() | Out-Null
